$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3601.7778
$ws.Range("I40").Value = 3230.8572
$ws.Range("J40").Value = 4900
$ws.Range("K40").Value = 3230.8572
$ws.Range("L40").Value = 4900
$ws.Range("M40").Value = -3055.8572
$ws.Range("N40").Value = -5250
$ws.Range("H43").Value = 161820.61
$ws.Range("I43").Value = 4493.1665
$ws.Range("K43").Value = 4493.1665
$ws.Range("M43").Value = -4424.1665
$ws.Range("H46").Value = 2835.5
$ws.Range("I46").Value = 7900
$ws.Range("J46").Value = 1822.6
$ws.Range("K46").Value = 23700
$ws.Range("L46").Value = 5467.799999999999
$ws.Range("M46").Value = -23581
$ws.Range("N46").Value = -5705.799999999999
$ws.Range("H60").Value = 2835.5
$ws.Range("I60").Value = 7900
$ws.Range("J60").Value = 1822.6
$ws.Range("K60").Value = 23700
$ws.Range("L60").Value = 5467.799999999999
$ws.Range("M60").Value = -23216
$ws.Range("N60").Value = -6435.799999999999
$ws.Range("H69").Value = 19189.125
$ws.Range("I69").Value = 11928.286
$ws.Range("K69").Value = 35784.858
$ws.Range("M69").Value = -34910.858
$ws.Range("H72").Value = 19189.125
$ws.Range("I72").Value = 11928.286
$ws.Range("K72").Value = 107354.574
$ws.Range("M72").Value = -102986.574
$ws.Range("H100").Value = 2634.182
$ws.Range("I100").Value = 1831
$ws.Range("K100").Value = 1831
$ws.Range("M100").Value = -1290
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988
$ws.Range("H132").Value = 2900
$ws.Range("I132").Value = 3009.2632
$ws.Range("K132").Value = 9027.7896
$ws.Range("M132").Value = -6497.7896
$ws.Range("H137").Value = 5926.9
$ws.Range("I137").Value = 1831.75
$ws.Range("K137").Value = 5495.25
$ws.Range("M137").Value = -2945.25
$ws.Range("H138").Value = 3204.26
$ws.Range("I138").Value = 1196.1666
$ws.Range("J138").Value = 3478.0908
$ws.Range("K138").Value = 3588.4998
$ws.Range("L138").Value = 10434.2724
$ws.Range("M138").Value = 1551.5002
$ws.Range("N138").Value = -20714.2724

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1364.5161
$ws.Range("I2").Value = 909.2917
$ws.Range("K2").Value = 909.2917
$ws.Range("M2").Value = -796.2917
$ws.Range("H97").Value = 1587.1875
$ws.Range("I97").Value = 1523.6666
$ws.Range("J97").Value = 1777.75
$ws.Range("K97").Value = 1523.6666
$ws.Range("L97").Value = 1777.75
$ws.Range("M97").Value = -1027.6666
$ws.Range("N97").Value = -2769.75
$ws.Range("H106").Value = 37793.8
$ws.Range("J106").Value = 37793.8
$ws.Range("L106").Value = 37793.8
$ws.Range("N106").Value = -40317.8
$ws.Range("H116").Value = 1364.5161
$ws.Range("I116").Value = 909.2917
$ws.Range("K116").Value = 909.2917
$ws.Range("M116").Value = 1384.7083
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H132").Value = 9182.714
$ws.Range("I132").Value = 3208.6667
$ws.Range("K132").Value = 9626.000100000001
$ws.Range("M132").Value = -7096.000100000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1364.5161
$ws.Range("I3").Value = 909.2917
$ws.Range("K3").Value = 909.2917
$ws.Range("M3").Value = -795.2917
$ws.Range("H105").Value = 2455.5557
$ws.Range("I105").Value = 1366.6666
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 1366.6666
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = 380.3334
$ws.Range("N105").Value = -6494

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 446425.56
$ws.Range("I31").Value = 5463.6665
$ws.Range("J31").Value = 1174012.6
$ws.Range("K31").Value = 5463.6665
$ws.Range("L31").Value = 1174012.6
$ws.Range("M31").Value = -5168.6665
$ws.Range("N31").Value = -1174602.6
$ws.Range("H34").Value = 446425.56
$ws.Range("I34").Value = 5463.6665
$ws.Range("J34").Value = 1174012.6
$ws.Range("K34").Value = 5463.6665
$ws.Range("L34").Value = 1174012.6
$ws.Range("M34").Value = -5261.6665
$ws.Range("N34").Value = -1174416.6
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H86").Value = 3841
$ws.Range("I86").Value = 3841
$ws.Range("K86").Value = 3841
$ws.Range("M86").Value = -2718
$ws.Range("H89").Value = 3841
$ws.Range("I89").Value = 3841
$ws.Range("K89").Value = 19205
$ws.Range("M89").Value = -13589
$ws.Range("H132").Value = 2299.3333
$ws.Range("I132").Value = 2081.647
$ws.Range("K132").Value = 6244.941
$ws.Range("M132").Value = -3714.941
$ws.Range("H134").Value = 527932.7
$ws.Range("I134").Value = 527932.7
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 1583798.1
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1581263.1
$ws.Range("N134").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 16175
$ws.Range("J3").Value = 19000
$ws.Range("L3").Value = 57000
$ws.Range("N3").Value = -57224
$ws.Range("H49").Value = 2148.2
$ws.Range("I49").Value = 1734.3334
$ws.Range("J49").Value = 2424.111
$ws.Range("K49").Value = 5203.0002
$ws.Range("L49").Value = 7272.333
$ws.Range("M49").Value = -5047.0002
$ws.Range("N49").Value = -7584.333
$ws.Range("H54").Value = 16123.75
$ws.Range("J54").Value = 17497.5
$ws.Range("L54").Value = 52492.5
$ws.Range("N54").Value = -53610.5
$ws.Range("H75").Value = 266667730
$ws.Range("J75").Value = 166667470
$ws.Range("L75").Value = 500002410
$ws.Range("N75").Value = -500004406
$ws.Range("H78").Value = 266667730
$ws.Range("J78").Value = 166667470
$ws.Range("L78").Value = 1500007230
$ws.Range("N78").Value = -1500017214
$ws.Range("H93").Value = 5715.727
$ws.Range("J93").Value = 5437.3
$ws.Range("L93").Value = 16311.9
$ws.Range("N93").Value = -20055.9
$ws.Range("H137").Value = 10000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 2962.36
$ws.Range("I139").Value = 4853
$ws.Range("J139").Value = 2704.5454
$ws.Range("K139").Value = 14559
$ws.Range("L139").Value = 8113.6362
$ws.Range("M139").Value = -9419
$ws.Range("N139").Value = -18393.6362

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2849.25
$ws.Range("I97").Value = 3086.2856
$ws.Range("J97").Value = 1190
$ws.Range("K97").Value = 3086.2856
$ws.Range("L97").Value = 1190
$ws.Range("M97").Value = -2590.2856
$ws.Range("N97").Value = -2182
$ws.Range("H122").Value = 1207.375
$ws.Range("I122").Value = 958.2727
$ws.Range("J122").Value = 1755.4
$ws.Range("K122").Value = 2874.8181
$ws.Range("L122").Value = 5266.200000000001
$ws.Range("M122").Value = -424.8181
$ws.Range("N122").Value = -10166.2
$ws.Range("H132").Value = 47621570
$ws.Range("I132").Value = 58826196
$ws.Range("J132").Value = 1903.25
$ws.Range("K132").Value = 176478588
$ws.Range("L132").Value = 5709.75
$ws.Range("M132").Value = -176476058
$ws.Range("N132").Value = -10769.75

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 147257.72
$ws.Range("I7").Value = 1594.6666
$ws.Range("K7").Value = 1594.6666
$ws.Range("M7").Value = -1482.6666
$ws.Range("H20").Value = 24139.357
$ws.Range("I20").Value = 19333.334
$ws.Range("J20").Value = 25450.092
$ws.Range("K20").Value = 19333.334
$ws.Range("L20").Value = 25450.092
$ws.Range("M20").Value = -19107.334
$ws.Range("N20").Value = -25902.092
$ws.Range("H93").Value = 66675748
$ws.Range("I93").Value = 66675748
$ws.Range("K93").Value = 66675748
$ws.Range("M93").Value = -66674500
$ws.Range("H122").Value = 5246.711
$ws.Range("J122").Value = 5769.8
$ws.Range("L122").Value = 17309.4
$ws.Range("N122").Value = -22209.4
$ws.Range("H126").Value = 147257.72
$ws.Range("I126").Value = 1594.6666
$ws.Range("K126").Value = 4783.9998
$ws.Range("M126").Value = -2313.9998
$ws.Range("H132").Value = 36320.832
$ws.Range("I132").Value = 3312.7273
$ws.Range("J132").Value = 127093.125
$ws.Range("K132").Value = 9938.1819
$ws.Range("L132").Value = 381279.375
$ws.Range("M132").Value = -7408.1819
$ws.Range("N132").Value = -386339.375

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()
$ws.Range("H132").Value = 1659.1628
$ws.Range("J132").Value = 1690.5454
$ws.Range("L132").Value = 5071.6362
$ws.Range("N132").Value = -10131.6362
